# custom accuracy + data trim
# 1) Round the last data row (row 5) values to 2 decimal places ("custom accuracy"),
# 2) Delete the now-redundant trailing row 6,
# 3) Narrow column N (14) from width 7 to width 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Apply "custom accuracy" (round to 2 decimals) to row 5, columns B..AH ---
$lastCol = 34  # column AH
for ($col = 2; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(5, $col)
    $orig = $cell.Value2
    $cell.Value = $excel.WorksheetFunction.Round($orig, 2)
}

# --- 2) Delete row 6 (trailing data row no longer needed) ---
$ws.Rows.Item(6).Delete()

# --- 3) Narrow column N (14th column) from width 7 to width 6 ---
$ws.Columns.Item(14).ColumnWidth = 5.17
